# 17.13.1.1 Macroeconomic Dashboard - extend the data table with 2021-2023
# columns (O, P, Q) and fix a handful of formatting inconsistencies in the
# existing columns A-C and N, matching the source data revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

function Copy-Format([string]$srcAddr, [string]$dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# 1. Columns A:C become one uniform, narrower width (was 43.14 / 46.29 /
#    43.57 characters individually, now ~37.43 for all three).
# ---------------------------------------------------------------------
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 36.67

# ---------------------------------------------------------------------
# 2. Reset the stored selection back to the top-left of the sheet.
# ---------------------------------------------------------------------
$ws.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Row formatting fix-ups: a few cells in columns B/C that should use
#    the wrapping style (same as column A) did not.
# ---------------------------------------------------------------------
Copy-Format "A4" "C4"
Copy-Format "A4" "B5"
Copy-Format "A4" "C5"
Copy-Format "A4" "B6"
Copy-Format "A4" "C6"
Copy-Format "A4" "B9"

# Row 9 grows taller once B9 wraps text, matching rows 5-7.
$ws.Rows.Item(9).RowHeight = 24

# ---------------------------------------------------------------------
# 4. Correct a transcription error in M6 (104.1 -> 104.26).
# ---------------------------------------------------------------------
$ws.Range("M6").Value = 104.26

# ---------------------------------------------------------------------
# 5. Fill in the previously-empty 2020 (N) column and add new 2021-2023
#    (O, P, Q) columns of data.
# ---------------------------------------------------------------------

# Header row (year labels) - reuse the formatting already used by N3.
Copy-Format "N3" "O3"
Copy-Format "N3" "P3"
Copy-Format "N3" "Q3"
$ws.Range("O3").Value = 2021
$ws.Range("P3").Value = 2022
$ws.Range("Q3").Value = 2023

# Row 4
Copy-Format "M4" "N4"
Copy-Format "M4" "O4"
Copy-Format "M4" "P4"
Copy-Format "M4" "Q4"
$ws.Range("N4").Value = 92.9
$ws.Range("O4").Value = 105.5
$ws.Range("P4").Value = 109
$ws.Range("Q4").Value = 106.2

# Row 5
Copy-Format "M5" "N5"
Copy-Format "M5" "O5"
Copy-Format "M5" "P5"
Copy-Format "M5" "Q5"
$ws.Range("N5").Value = 106.3
$ws.Range("O5").Value = 111.90503981851454
$ws.Range("P5").Value = 113.92290931741762
$ws.Range("Q5").Value = 110.8

# Row 6
Copy-Format "M6" "N6"
Copy-Format "M6" "O6"
Copy-Format "M6" "P6"
Copy-Format "M6" "Q6"
$ws.Range("N6").Value = 121.27
$ws.Range("O6").Value = 111.5
$ws.Range("P6").Value = 105.1
$ws.Range("Q6").Value = 109.3

# Row 7 (Q7 stays blank, just gets the formatting)
Copy-Format "M7" "N7"
Copy-Format "M7" "O7"
Copy-Format "M7" "P7"
Copy-Format "M7" "Q7"
$ws.Range("N7").Value = -19734.0366
$ws.Range("O7").Value = -1763.6
$ws.Range("P7").Value = -10400.700000000001

# Row 8 (M8 was blank before, Q8 stays blank)
Copy-Format "L8" "M8"
Copy-Format "L8" "N8"
Copy-Format "L8" "O8"
Copy-Format "L8" "P8"
Copy-Format "L8" "Q8"
$ws.Range("M8").Value = 319474.59999999998
$ws.Range("N8").Value = 407116.85000000003
$ws.Range("O8").Value = 436586.8
$ws.Range("P8").Value = 477967.8

# Row 9 (Q9 stays blank)
Copy-Format "M9" "N9"
Copy-Format "M9" "O9"
Copy-Format "M9" "P9"
Copy-Format "M9" "Q9"
$ws.Range("N9").Value = 81.599999999999994
$ws.Range("O9").Value = 146.4
$ws.Range("P9").Value = 144.69999999999999

# Row 10 (Q10 stays blank)
Copy-Format "M10" "N10"
Copy-Format "M10" "O10"
Copy-Format "M10" "P10"
Copy-Format "M10" "Q10"
$ws.Range("N10").Value = 1973.2
$ws.Range("O10").Value = 2752.1
$ws.Range("P10").Value = 2254.6999999999998

# Row 11 (Q11 stays blank)
Copy-Format "M11" "N11"
Copy-Format "M11" "O11"
Copy-Format "M11" "P11"
Copy-Format "M11" "Q11"
$ws.Range("N11").Value = 3718.8
$ws.Range("O11").Value = 5580.2
$ws.Range("P11").Value = 9803.2000000000007

$ws.Range("A1").Select()
